$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$find,
        [string]$replace,
        [switch]$MatchCase
    )
    $r = $d.Content
    $matchCaseFlag = if ($MatchCase) { $true } else { $false }
    $ok = $r.Find.Execute($find, $matchCaseFlag, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "FAILED to find: $find"
    }
    return $ok
}

function Delete-ParagraphExact {
    param(
        [string]$exactText
    )
    $found = $false
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text
        # Trim the trailing paragraph mark / cell mark character(s) for comparison.
        $trimmed = $t.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $exactText) {
            $p.Range.Delete()
            $found = $true
            break
        }
    }
    if (-not $found) {
        Write-Output "FAILED to find paragraph with exact text: $exactText"
    }
    return $found
}

# 1. Team member table: merge "Ee" + " Lyn Lim" into a single run (also removes the
#    _GoBack bookmark that used to sit between them -- it is re-created later, near
#    the end of the document).
Replace-Text "Ee Lyn Lim" "Ee Lyn Lim"

# 2. "Things that went well" bullet list
Replace-Text "Everyone collaborated well together and respectful of each other" `
             "Everyone collaborated well together and are respectful of each other"

Replace-Text "Tasks delegated are carried out by relevant team members (showing responsibility)" `
             "Tasks delegated are carried out by relevant team members "

Replace-Text "Team members are willing to help others when required" `
             "Team members are willing to help each other when required"

Replace-Text "Team members do not hesitate to get feedback from each other or the tutor when a problem or uncertainty arises" `
             "Most team members do not hesitate to get feedback from each other or the tutor when a problem or uncertainty arises"

Replace-Text "Team members incorporates advice and feedback into their work" `
             "Team members incorporate advice and feedback into their work"

Replace-Text "Source control (i.e. backing up work regularly)" `
             "Team has good source control (i.e. backing up work regularly)"

# 3. "Things that could have gone better" bullet list
Replace-Text "Learning Django framework (researching and getting information on Django, Python, MySQL integration syntax)" `
             "Learning Django framework (researching and getting information on Django, Python, MySQL integration syntax) by all team members"

Replace-Text "More participation or brainstorming of ideas from all team members" `
             "More motivation, initiative and participation from all team members"

Replace-Text "Correct story points to time estimation depiction" `
             "Making sure all data presented are accurate (i.e. story points to time estimation depiction)"

Replace-Text "More communication between developer and client (such as, gaining detailed specifications from the client before development – would help better understand data flow)" `
             "Regular communication between developer and client "

Replace-Text "Ensure impartial division of tasks between team members" `
             "Ensure fair division of tasks between team members and that deadlines are met"

# 4. "Things that surprised us" bullet list
Replace-Text "Learning syntaxes and debugging took longer than expected" `
             "Debugging and editing took longer than expected "

Replace-Text "Django framework and database management were more complex than expected" `
             "Learning the syntaxes to Django template system and integrating the database were challenging, resulting to demotivation among most team members"

Replace-Text "Database integration proved to be challenging" `
             "Due to different priorities among team members, it was difficult to ensure every team member puts in equal amount of time and effort into the project "

Replace-Text "Editing and drafting of work was time consuming" `
             "Unaware that Django had automated features available (i.e. authentication), which could have made the project development a lot less complicated and less time consuming"

Delete-ParagraphExact "Unaware that Django had automated features available (i.e. authentication)"

# 5. "Ways to improve for the next Sprint" bullet list
Replace-Text "estimation for each story point, making sure that each story point accounts to 8 hours" `
             "estimation for each story point, so that each story point accounts for 8 hours"

# 6. Summary section
Replace-Text "To ensure effective efforts and quality output, the team will continue to maintain transparent communication and project progress. Also, to continuously seek feedback from the tutor. " `
             "The team will continue to maintain transparent communication and project progress. Simultaneously, the team will continue to seek feedback from the tutor, check specifications and generate test cases as a method of quality control. "

Replace-Text "itemise Sprint tasks to be more in-depth" `
             "itemise Sprint tasks to be more in-depth and to delegate work more effectively,"

Replace-Text "delegate work more effectively,    " `
             "put in equal effort in learning Django framework and utilise features available, "

Replace-Text "carefully draft work to ensure accurate data is presented (i.e. story points, time estimation)," `
             "carefully draft work to ensure accurate data are presented (i.e. story points, time estimation) "

Delete-ParagraphExact "establishing quality control by checking specifications and using test cases"

# Re-insert the _GoBack bookmark at its new location, inside the "carefully draft..."
# bullet, right after "accurate data are".
$r2 = $d.Content
if ($r2.Find.Execute("accurate data are", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $bm = $d.Range($r2.End, $r2.End)
    $d.Bookmarks.Add("_GoBack", $bm) | Out-Null
}
